# Populate a CFDI-style "Conceptos" table on Sheet1 with a bold, bordered,
# centered header row and the corresponding data row below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("Cantidad", "ClaveProdServ", "ClaveUnidad", "Descripcion", "Importe", "NoIdentificacion", "ObjetoImp", "Unidad", "ValorUnitario")
$data    = @("77.58", "15101505", "LTR", "Diesel (Despacho 4024741-0)", "1729.20", "PL/8748/EXP/ES/2015-4024741", "02", "Litros", "22.28925")

# Columns that hold numeric-looking text which must stay as text (not be
# coerced into a number by Excel's automatic type detection).
$textColumns = @(1, 2, 5, 7, 9)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

for ($i = 0; $i -lt $data.Length; $i++) {
    $cell = $ws.Cells.Item(2, $i + 1)
    if ($textColumns -contains ($i + 1)) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $data[$i]
}

$headerRange = $ws.Range("A1:I1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = "thin"
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
